# Update header labels (row 1) so that Power BI can automatically turn
# the first row into a header when the table is loaded.
# Year-only headers get prefixed with "Ano " and the interval-based
# sheet gets its headers prefixed with "Intervalo ".

$wb = $excel.ActiveWorkbook

# Sheets 1-3 and 5: simple "Ano <year>" headers in B1:E1
$anoSheets = @(
    "Potencia Acumulada - SIN (MW)",
    "Geracao Periodo Medio (MWMed)",
    "Atendimento a Ponta(MW)",
    "Emissoes Totais (MtCO2eq)"
)

foreach ($sheetName in $anoSheets) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("B1").Value = "Ano " + $ws.Range("B1").Value2
    $ws.Range("C1").Value = "Ano " + $ws.Range("C1").Value2
    $ws.Range("D1").Value = "Ano " + $ws.Range("D1").Value2
    $ws.Range("E1").Value = "Ano " + $ws.Range("E1").Value2
}

# Sheet 4: "Intervalo <period>" headers in B1:E1
$ws4 = $wb.Worksheets.Item("Potencia Incremental - SIN(MW)")
$ws4.Range("B1").Value = "Intervalo " + $ws4.Range("B1").Value2
$ws4.Range("C1").Value = "Intervalo " + $ws4.Range("C1").Value2
$ws4.Range("D1").Value = "Intervalo " + $ws4.Range("D1").Value2
$ws4.Range("E1").Value = "Intervalo " + $ws4.Range("E1").Value2

# Sheet 6: only B1 header present, "Ano <year>"
$ws6 = $wb.Worksheets.Item("Custo Total (bilhões de R$)")
$ws6.Range("B1").Value = "Ano " + $ws6.Range("B1").Value2
